# Re-order / refresh the "line file" rows used for the highcharts
# visualization: the row order changes (Line:8 rows move up, a new
# Line:3 block of rows is appended, Line:9's "Pri pH flows" block
# shifts down) and every row's uuid is refreshed to the new batch id
# (2b1ef34e-1545-47d2-9511-6f993db986ca). Sheet grows from 15 to 19
# data-bearing rows (A1:H15 -> A1:H19).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "Line:8 Stage:1"
$ws.Cells.Item(2,2).Value = "'01/09/2024"
$ws.Cells.Item(2,3).Value = "pri cl LA"
$ws.Cells.Item(2,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2,4).Value = 45300.42056299769
$ws.Cells.Item(2,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2,5).Value = 45300.42067873842
$ws.Cells.Item(2,6).Value = 0.17
$ws.Cells.Item(2,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(2,8).NumberFormat = "0"
$ws.Cells.Item(2,8).Value = 0.0001157407407407407

$ws.Cells.Item(3,1).Value = "Line:8 Stage:1"
$ws.Cells.Item(3,2).Value = "'01/09/2024"
$ws.Cells.Item(3,3).Value = "pri pH HA"
$ws.Cells.Item(3,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3,4).Value = 45300.65476473379
$ws.Cells.Item(3,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3,5).Value = 45300.65488047454
$ws.Cells.Item(3,6).Value = 0.17
$ws.Cells.Item(3,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(3,8).NumberFormat = "0"
$ws.Cells.Item(3,8).Value = 0.0001157407407407407

$ws.Cells.Item(4,1).Value = "Line:8 Stage:1"
$ws.Cells.Item(4,2).Value = "'01/09/2024"
$ws.Cells.Item(4,3).Value = "pri cl LA"
$ws.Cells.Item(4,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4,4).Value = 45300.40887037037
$ws.Cells.Item(4,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4,5).Value = 45300.40910185185
$ws.Cells.Item(4,6).Value = 0.33
$ws.Cells.Item(4,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(4,8).NumberFormat = "0"
$ws.Cells.Item(4,8).Value = 0.0002314814814814815

$ws.Cells.Item(5,1).Value = "Line:9 Stage:1"
$ws.Cells.Item(5,2).Value = "'01/09/2024"
$ws.Cells.Item(5,3).Value = "Pri pH flows"
$ws.Cells.Item(5,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5,4).Value = 45300.49280315972
$ws.Cells.Item(5,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5,5).Value = 45300.49893741898
$ws.Cells.Item(5,6).Value = 8.83
$ws.Cells.Item(5,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(5,8).NumberFormat = "0"
$ws.Cells.Item(5,8).Value = 0.006134259259259259

$ws.Cells.Item(6,1).Value = "Line:9 Stage:1"
$ws.Cells.Item(6,2).Value = "'01/09/2024"
$ws.Cells.Item(6,3).Value = "Pri pH flows"
$ws.Cells.Item(6,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6,4).Value = 45300.50241207176
$ws.Cells.Item(6,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6,5).Value = 45300.50287503472
$ws.Cells.Item(6,6).Value = 0.67
$ws.Cells.Item(6,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(6,8).NumberFormat = "0"
$ws.Cells.Item(6,8).Value = 0.000462962962962963

$ws.Cells.Item(7,1).Value = "Line:9 Stage:1"
$ws.Cells.Item(7,2).Value = "'01/09/2024"
$ws.Cells.Item(7,3).Value = "Pri pH flows"
$ws.Cells.Item(7,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7,4).Value = 45300.50356979167
$ws.Cells.Item(7,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7,5).Value = 45300.50368553241
$ws.Cells.Item(7,6).Value = 0.17
$ws.Cells.Item(7,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(7,8).NumberFormat = "0"
$ws.Cells.Item(7,8).Value = 0.0001157407407407407

$ws.Cells.Item(8,1).Value = "Line:9 Stage:1"
$ws.Cells.Item(8,2).Value = "'01/09/2024"
$ws.Cells.Item(8,3).Value = "Pri pH flows"
$ws.Cells.Item(8,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8,4).Value = 45300.50438012731
$ws.Cells.Item(8,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8,5).Value = 45300.50449586806
$ws.Cells.Item(8,6).Value = 0.17
$ws.Cells.Item(8,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(8,8).NumberFormat = "0"
$ws.Cells.Item(8,8).Value = 0.0001157407407407407

$ws.Cells.Item(9,1).Value = "Line:9 Stage:1"
$ws.Cells.Item(9,2).Value = "'01/09/2024"
$ws.Cells.Item(9,3).Value = "Pri pH flows"
$ws.Cells.Item(9,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9,4).Value = 45300.50519054398
$ws.Cells.Item(9,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9,5).Value = 45300.50715813658
$ws.Cells.Item(9,6).Value = 2.83
$ws.Cells.Item(9,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(9,8).NumberFormat = "0"
$ws.Cells.Item(9,8).Value = 0.001967592592592592

$ws.Cells.Item(10,1).Value = "Line:3 Stage:1"
$ws.Cells.Item(10,2).Value = "'12/11/2023"
$ws.Cells.Item(10,3).Value = "SEC Cl"
$ws.Cells.Item(10,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10,4).Value = 45271.36145825232
$ws.Cells.Item(10,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10,5).Value = 45271.36158556713
$ws.Cells.Item(10,6).Value = 0.18
$ws.Cells.Item(10,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(10,8).NumberFormat = "0"
$ws.Cells.Item(10,8).Value = 0.0001273148148148148

$ws.Cells.Item(11,1).Value = "Line:3 Stage:1"
$ws.Cells.Item(11,2).Value = "'12/11/2023"
$ws.Cells.Item(11,3).Value = "PRI pH"
$ws.Cells.Item(11,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11,4).Value = 45271.68537890046
$ws.Cells.Item(11,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11,5).Value = 45271.68549464121
$ws.Cells.Item(11,6).Value = 0.17
$ws.Cells.Item(11,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(11,8).NumberFormat = "0"
$ws.Cells.Item(11,8).Value = 0.0001157407407407407

$ws.Cells.Item(12,1).Value = "Line:3 Stage:1"
$ws.Cells.Item(12,2).Value = "'12/11/2023"
$ws.Cells.Item(12,3).Value = "PRI pH SEC pH"
$ws.Cells.Item(12,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12,4).Value = 45271.68549475694
$ws.Cells.Item(12,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12,5).Value = 45271.68769383102
$ws.Cells.Item(12,6).Value = 3.17
$ws.Cells.Item(12,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(12,8).NumberFormat = "0"
$ws.Cells.Item(12,8).Value = 0.002199074074074074

$ws.Cells.Item(13,1).Value = "Line:3 Stage:1"
$ws.Cells.Item(13,2).Value = "'12/11/2023"
$ws.Cells.Item(13,3).Value = "PRI pH&rem SEC pH"
$ws.Cells.Item(13,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13,4).Value = 45271.6876965625
$ws.Cells.Item(13,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13,5).Value = 45271.68898128472
$ws.Cells.Item(13,6).Value = 1.85
$ws.Cells.Item(13,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(13,8).NumberFormat = "0"
$ws.Cells.Item(13,8).Value = 0.001284722222222222

$ws.Cells.Item(14,1).Value = "Line:3 Stage:1"
$ws.Cells.Item(14,2).Value = "'12/11/2023"
$ws.Cells.Item(14,3).Value = "PRI pH&rem SEC pH&rem"
$ws.Cells.Item(14,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14,4).Value = 45271.68897129629
$ws.Cells.Item(14,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14,5).Value = 45271.68908703703
$ws.Cells.Item(14,6).Value = 0.17
$ws.Cells.Item(14,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(14,8).NumberFormat = "0"
$ws.Cells.Item(14,8).Value = 0.0001157407407407407

$ws.Cells.Item(15,1).Value = "Line:3 Stage:1"
$ws.Cells.Item(15,2).Value = "'12/11/2023"
$ws.Cells.Item(15,3).Value = "PRI pH SEC pH"
$ws.Cells.Item(15,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15,4).Value = 45271.6890871875
$ws.Cells.Item(15,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15,5).Value = 45271.68920292824
$ws.Cells.Item(15,6).Value = 0.17
$ws.Cells.Item(15,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(15,8).NumberFormat = "0"
$ws.Cells.Item(15,8).Value = 0.0001157407407407407

$ws.Cells.Item(16,1).Value = "Line:3 Stage:1"
$ws.Cells.Item(16,2).Value = "'12/11/2023"
$ws.Cells.Item(16,3).Value = "PRI pH SEC pH"
$ws.Cells.Item(16,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16,4).Value = 45271.89825825232
$ws.Cells.Item(16,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16,5).Value = 45271.8996471412
$ws.Cells.Item(16,6).Value = 2
$ws.Cells.Item(16,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(16,8).NumberFormat = "0"
$ws.Cells.Item(16,8).Value = 0.001388888888888889

$ws.Cells.Item(17,1).Value = "Line:3 Stage:1"
$ws.Cells.Item(17,2).Value = "'12/11/2023"
$ws.Cells.Item(17,3).Value = "PRI Cl&pH"
$ws.Cells.Item(17,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17,4).Value = 45272.09804105324
$ws.Cells.Item(17,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17,5).Value = 45272.09827253472
$ws.Cells.Item(17,6).Value = 0.33
$ws.Cells.Item(17,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(17,8).NumberFormat = "0"
$ws.Cells.Item(17,8).Value = 0.0002314814814814815

$ws.Cells.Item(18,1).Value = "Line:8 Stage:1"
$ws.Cells.Item(18,2).Value = "'01/09/2024"
$ws.Cells.Item(18,3).Value = "pri cl LA"
$ws.Cells.Item(18,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18,4).Value = 45300.42056299769
$ws.Cells.Item(18,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18,5).Value = 45300.42067873842
$ws.Cells.Item(18,6).Value = 0.17
$ws.Cells.Item(18,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(18,8).NumberFormat = "0"
$ws.Cells.Item(18,8).Value = 0.0001157407407407407

$ws.Cells.Item(19,1).Value = "Line:8 Stage:1"
$ws.Cells.Item(19,2).Value = "'01/09/2024"
$ws.Cells.Item(19,3).Value = "pri pH HA"
$ws.Cells.Item(19,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(19,4).Value = 45300.65476473379
$ws.Cells.Item(19,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(19,5).Value = 45300.65488047454
$ws.Cells.Item(19,6).Value = 0.17
$ws.Cells.Item(19,7).Value = "2b1ef34e-1545-47d2-9511-6f993db986ca"
$ws.Cells.Item(19,8).NumberFormat = "0"
$ws.Cells.Item(19,8).Value = 0.0001157407407407407

